$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the width of the new "IP" column (K) to match the target (15.5 chars)
$ws.Columns("K").ColumnWidth = 14.625

# Add the new "IP" column header and values.
# Order matters: it determines the order new entries are appended to the
# shared string table, matching the author's original entry order
# (header first, then row 7 / iMac #6, then row 6 / iMac #5, then row 8 / iMac #7).
$ws.Range("K1").Value = "IP"
$ws.Range("K7").Value = "205.208.43.103"
$ws.Range("K6").Value = "205.208.92.231"
$ws.Range("K8").Value = "205.208.43.70"

# Match the final selection left by the author
$ws.Range("K9").Select() | Out-Null
